# Apply updated cryptocurrency price/volume data to Sheet1 (columns D and E, rows 2-51).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some "Price" values are plain decimals (e.g. "535.12") that Excel would
# otherwise auto-convert to numbers; force those specific cells to Text format
# first so they stay text, matching the source data (prices with "." as a
# thousands separator, like "58.553.94", remain text automatically).
$textCells = @("D5","D6","D8","D16","D20","D23","D24","D25","D28","D32","D33","D34","D35","D36","D37","D38","D39","D42","D44","D46","D47","D48")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '58.553.94'
$ws.Range("E2").Value = '  -1.82%  '
$ws.Range("D3").Value = '2.615.05'
$ws.Range("E3").Value = '  +0.01%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '535.12'
$ws.Range("E5").Value = '  -0.60%  '
$ws.Range("D6").Value = '142.78'
$ws.Range("E6").Value = '  +0.38%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").Value = '0.586'
$ws.Range("E8").Value = '  +3.40%  '
$ws.Range("D9").Value = '2.617.66'
$ws.Range("E9").Value = '  -0.19%  '
$ws.Range("E10").Value = '  +3.85%  '
$ws.Range("E11").Value = '  -2.07%  '
$ws.Range("E12").Value = '  -0.71%  '
$ws.Range("E13").Value = '  +2.12%  '
$ws.Range("D14").Value = '3.069.38'
$ws.Range("E14").Value = '  -0.18%  '
$ws.Range("D15").Value = '58.520.86'
$ws.Range("E15").Value = '  -1.79%  '
$ws.Range("D16").Value = '20.81'
$ws.Range("E16").Value = '  +0.17%  '
$ws.Range("D17").Value = '2.617.11'
$ws.Range("E17").Value = '  -1.65%  '
$ws.Range("E18").Value = '  -1.43%  '
$ws.Range("E19").Value = '  +1.08%  '
$ws.Range("D20").Value = '334.56'
$ws.Range("E20").Value = '  -1.85%  '
$ws.Range("E21").Value = '  +0.00%  '
$ws.Range("E22").Value = '  -2.66%  '
$ws.Range("D23").Value = '0.998'
$ws.Range("E23").Value = '  -0.11%  '
$ws.Range("D24").Value = '66.91'
$ws.Range("E24").Value = '  -0.56%  '
$ws.Range("D25").Value = '0.420'
$ws.Range("E25").Value = '  +2.66%  '
$ws.Range("E26").Value = '  -0.13%  '
$ws.Range("E27").Value = '  -3.02%  '
$ws.Range("D28").Value = '7.08'
$ws.Range("E28").Value = '  -2.17%  '
$ws.Range("E29").Value = '  -1.69%  '
$ws.Range("E30").Value = '  +0.02%  '
$ws.Range("E31").Value = '  -1.53%  '
$ws.Range("D32").Value = '5.92'
$ws.Range("E32").Value = '  +1.34%  '
$ws.Range("D33").Value = '153.27'
$ws.Range("E33").Value = '  +1.63%  '
$ws.Range("D34").Value = '18.89'
$ws.Range("E34").Value = '  +0.26%  '
$ws.Range("D35").Value = '3.90'
$ws.Range("E35").Value = '  -2.33%  '
$ws.Range("D36").Value = '37.11'
$ws.Range("E36").Value = '  -0.98%  '
$ws.Range("D37").Value = '1.11'
$ws.Range("E37").Value = '  -1.91%  '
$ws.Range("D38").Value = '0.835'
$ws.Range("E38").Value = '  -0.02%  '
$ws.Range("D39").Value = '0.825'
$ws.Range("E39").Value = '  -0.50%  '
$ws.Range("E40").Value = '  -2.86%  '
$ws.Range("E41").Value = '  +1.33%  '
$ws.Range("D42").Value = '283.97'
$ws.Range("E42").Value = '  +2.14%  '
$ws.Range("E43").Value = '  +0.01%  '
$ws.Range("D44").Value = '0.595'
$ws.Range("E44").Value = '  -1.33%  '
$ws.Range("E45").Value = '  -0.47%  '
$ws.Range("D46").Value = '0.0950'
$ws.Range("E46").Value = '  -0.09%  '
$ws.Range("D47").Value = '19.04'
$ws.Range("E47").Value = '  +1.25%  '
$ws.Range("D48").Value = '0.0530'
$ws.Range("E48").Value = '  +1.12%  '
$ws.Range("E49").Value = '  +0.92%  '
$ws.Range("D50").Value = '1.939.90'
$ws.Range("E50").Value = '  -0.44%  '
$ws.Range("E51").Value = '  -1.29%  '
